$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 8 (which holds "extr1"),
# shifting the existing "extr1".."extr8" rows down to rows 10-17.
$ws.Rows.Item(8).Resize(2).Insert()

# New row 8: line7
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# New row 9: line8
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Apply the same formatting as the other "A" column entries to the new rows
# (bold centered font with a thin box border — matches style index 1).
$rngA = $ws.Range("A8:A9")
$rngA.Font.Bold = $true
$rngA.HorizontalAlignment = -4108
$rngA.VerticalAlignment = -4160
$rngA.Borders.LineStyle = 1

# Update shifted rows (now 10-17, formerly extr1..extr8 at rows 8-15) with new data values.
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $false

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $false

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $false

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $false

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $false
